# Updated results with infer-no-filter
$wb = $excel.ActiveWorkbook

$wsAllTools = $wb.Worksheets.Item("all_tools")
$wsInfer    = $wb.Worksheets.Item("infer")

# ---------------------------------------------------------------------------
# Sheet "all_tools" rows 10-12 (tool id 6, dataset "time")
# ---------------------------------------------------------------------------
$wsAllTools.Range("G10").Value = 813
$wsAllTools.Range("I10").Value = -0.07157910106056362
$wsAllTools.Range("J10").Value = 0.4916600219100405
$wsAllTools.Range("K10").Value = -0.08830673038161191
$wsAllTools.Range("L10").Value = 0.5419792899048017

$wsAllTools.Range("G11").Value = 813
$wsAllTools.Range("I11").Value = -0.05507056613029693
$wsAllTools.Range("J11").Value = 0.5800104076897017
$wsAllTools.Range("K11").Value = -0.07444201065306216
$wsAllTools.Range("L11").Value = 0.6074026216973724

$wsAllTools.Range("G12").Value = 813
$wsAllTools.Range("I12").Value = 0.1102028102074909
$wsAllTools.Range("J12").Value = 0.265108023071319
$wsAllTools.Range("K12").Value = 0.1660257192865707
$wsAllTools.Range("L12").Value = 0.2491959671429019

# ---------------------------------------------------------------------------
# Sheet "all_tools" rows 25-29 (tool id "f", dataset "readability_level_ba")
# ---------------------------------------------------------------------------
$wsAllTools.Range("G25").Value = 39
$wsAllTools.Range("I25").Value = -0.1807753815155468
$wsAllTools.Range("J25").Value = 0.3541954904764164
$wsAllTools.Range("K25").Value = -0.2576049186596542
$wsAllTools.Range("L25").Value = 0.3354345184285685

$wsAllTools.Range("G26").Value = 39
$wsAllTools.Range("I26").Value = -0.1807753815155468
$wsAllTools.Range("J26").Value = 0.3541954904764164
$wsAllTools.Range("K26").Value = -0.2666436877354316
$wsAllTools.Range("L26").Value = 0.3181414648703181

$wsAllTools.Range("G27").Value = 39
$wsAllTools.Range("I27").Value = 0.3539900381483285
$wsAllTools.Range("J27").Value = 0.07056136851892029
$wsAllTools.Range("K27").Value = 0.4341802833034056
$wsAllTools.Range("L27").Value = 0.09288178063084394

$wsAllTools.Range("G28").Value = 39
$wsAllTools.Range("K28").Value = -0.2493004677260264
$wsAllTools.Range("L28").Value = 0.3517858440384553

$wsAllTools.Range("G29").Value = 39
$wsAllTools.Range("K29").Value = -0.1491396897503261
$wsAllTools.Range("L29").Value = 0.5814513259975999

# ---------------------------------------------------------------------------
# Sheet "infer" rows 10-12 (tool id 6, dataset "time")
# ---------------------------------------------------------------------------
$wsInfer.Range("F10").Value = 23
$wsInfer.Range("G10").Value = 24
$wsInfer.Range("I10").Value = -0.1454025530693833
$wsInfer.Range("J10").Value = 0.2372373518450496
$wsInfer.Range("K10").Value = -0.17271903862684
$wsInfer.Range("L10").Value = 0.2303502122764337

$wsInfer.Range("F11").Value = 23
$wsInfer.Range("G11").Value = 24
$wsInfer.Range("I11").Value = -0.1395616700784287
$wsInfer.Range("J11").Value = 0.2348980869048207
$wsInfer.Range("K11").Value = -0.1674579385094694
$wsInfer.Range("L11").Value = 0.2450782275649824

$wsInfer.Range("F12").Value = 23
$wsInfer.Range("G12").Value = 24
$wsInfer.Range("I12").Value = 0.02140819589682411
$wsInfer.Range("J12").Value = 0.8544862615484419
$wsInfer.Range("K12").Value = 0.02708713119452734
$wsInfer.Range("L12").Value = 0.8518765230635053

# ---------------------------------------------------------------------------
# Sheet "infer" rows 25-29 (tool id "f", dataset "readability_level_ba")
# Previously these rows had no correlation values (F/G were 0); now they do.
# ---------------------------------------------------------------------------
$wsInfer.Range("F25").Value = 1
$wsInfer.Range("G25").Value = 1
$wsInfer.Range("I25").Value = -0.3535533905932737
$wsInfer.Range("J25").Value = 0.1037416782365415
$wsInfer.Range("K25").Value = -0.4200840252084029
$wsInfer.Range("L25").Value = 0.105228057983522

$wsInfer.Range("F26").Value = 1
$wsInfer.Range("G26").Value = 1
$wsInfer.Range("I26").Value = -0.1649915822768611
$wsInfer.Range("J26").Value = 0.4476990724652935
$wsInfer.Range("K26").Value = -0.1960392117639214
$wsInfer.Range("L26").Value = 0.4668248490265503

$wsInfer.Range("F27").Value = 1
$wsInfer.Range("G27").Value = 1
$wsInfer.Range("I27").Value = 0.02366905341655754
$wsInfer.Range("J27").Value = 0.9135633303377861
$wsInfer.Range("K27").Value = 0.02802621677476181
$wsInfer.Range("L27").Value = 0.9179387985999929

$wsInfer.Range("F28").Value = 1
$wsInfer.Range("G28").Value = 1
$wsInfer.Range("I28").Value = -0.2625754538144587
$wsInfer.Range("J28").Value = 0.2314460271038938
$wsInfer.Range("K28").Value = -0.3089716991054783
$wsInfer.Range("L28").Value = 0.2442606266224961

$wsInfer.Range("F29").Value = 1
$wsInfer.Range("G29").Value = 1
$wsInfer.Range("I29").Value = 0.2592724864350675
$wsInfer.Range("J29").Value = 0.2328233516916538
$wsInfer.Range("K29").Value = 0.3080616184861621
$wsInfer.Range("L29").Value = 0.2457251662216493
